# Insert a new weekly price record at the top of the Acelga / Macroferia Regional
# de Talca data block (row 249), shifting all existing rows (249-286) down by one
# (to 250-287), and fill the new row 249 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 249; this shifts rows 249:286 down to 250:287
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row 249 with the new data record
$ws.Range("A249").Value = 5
$ws.Range("B249").Value = "Macroferia Regional de Talca"
$ws.Range("C249").Value = "Maule"
$ws.Range("D249").Value = 44776
$ws.Range("E249").Value = 7
$ws.Range("F249").Value = 100112009
$ws.Range("G249").Value = "Acelga"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 400
$ws.Range("K249").Value = 3500
$ws.Range("L249").Value = 3500
$ws.Range("M249").Value = 3500
$ws.Range("N249").Value = "$/docena de atados (4 kilos)"
$ws.Range("O249").Value = "Región del Maule"
$ws.Range("P249").Value = 875
$ws.Range("Q249").Value = 4
$ws.Range("R249").Value = "Hortaliza"

Write-Host "Row inserted and populated."
